# Apply "Minor improvements to setup, meta, forest, FAQs" edits to the
# network-to-do workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("To do")

# --- Row 46: setup / augment note -> mark as improved -----------------
$ws.Range("E46").Value = "improved"
$ws.Range("F46").Value = 43194

# --- Row 70-72: FAQ items marked done ----------------------------------
$ws.Range("E70").Value = "done"
$ws.Range("F70").Value = 43194

$ws.Range("E71").Value = "done"
$ws.Range("F71").Value = 42207

$ws.Range("E72").Value = "done"
$ws.Range("F72").Value = 42380

# --- Row 77, 79: mark done with completion dates -----------------------
$ws.Range("E77").Value = "done"
$ws.Range("F77").Value = 43194

$ws.Range("E79").Value = "done"
$ws.Range("F79").Value = 42207

# --- Row 80: too late to action -----------------------------------------
$ws.Range("E80").Value = "too late to be worth doing"
$ws.Range("F80").Value = "NA"

# --- Row 73: add a note, and grow the row to fit the extra line -------
$ws.Range("E73").Value = "not sure how to do this without making labels too long."
$ws.Rows.Item(73).RowHeight = 45

# --- Row 81: remove the underscore-naming-advice task (superseded) -----
$ws.Range("B81").FormatConditions().Item(1).Delete()
$ws.Rows.Item(81).Delete()

# --- Fix up the autofilter / filter-database range now that the sheet
#     only runs to row 80 -------------------------------------------------
$ws.AutoFilterMode = $false
$ws.Range("A2:F80").AutoFilter() | Out-Null
$fdb = $wb.Names().Item("_xlnm._FilterDatabase")
$fdb.RefersTo = "='To do'!`$A`$2:`$F`$80"

# --- Update the view so the frozen pane shows the newly edited rows ----
$ws.Activate()
$ws.Range("C26").Select()
